# ---------------------------------------------------------------------------
# Transform before.docx (generic "custom-reference" style showcase) into the
# "working with using custom template" document:
#   - Title/Subtitle/Author/Date paragraphs become Heading1..Heading4 samples
#   - Everything else in the middle of the document (old Heading1-9 samples,
#     FirstParagraph, BodyText w/ hyperlink+footnote, BlockText, TableCaption,
#     the demo table, ImageCaption, DefinitionTerm/Definition pairs) is removed
#   - New trailing content is appended: a blank paragraph, a "Normal text "
#     paragraph, another blank paragraph, a List Bullet paragraph and a
#     List Number paragraph.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# 1) Remove the footnote definition(s) up front (the only footnote reference
#    lives inside the "Body Text" paragraph that we are about to delete).
while ($d.Footnotes.Count -gt 0) {
    $d.Footnotes.Item(1).Delete()
}

# 2) Delete the trailing paragraphs that come *after* the demo table:
#    Image Caption, DefinitionTerm, Definition, DefinitionTerm, Definition.
#    (indices 25..29 in the original document) -- delete from the end so
#    earlier indices remain valid.
for ($i = 29; $i -ge 25; $i--) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Delete()
}

# 3) Delete the paragraphs between "Date" (kept) and the table:
#    Abstract, the old Heading1..Heading9 bookmark samples, FirstParagraph,
#    BodyText, BlockText, TableCaption (indices 5..18). Again delete from the
#    end backwards so we never invalidate the indices we still need.
for ($i = 18; $i -ge 5; $i--) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Delete()
}

# 4) Turn the 4 remaining lead paragraphs into the new Heading 1..4 samples.
#    (Must happen *before* the table is removed -- deleting the table makes
#    further indexed Paragraphs.Item(...).Range reads/writes unreliable.)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Heading 1 Style"
$p1.Style = "Heading 1"

$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Heading 2 Style"
$p2.Style = "Heading 2"

$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "Heading 3 Style"
$p3.Style = "Heading 3"

$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "Heading 4 Style"
$p4.Style = "Heading 4"

# 5) Append the new trailing paragraphs (still before removing the table, for
#    the same reliability reason as above).

# 5a) blank paragraph
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Style = "Normal"

# 5b) "Normal text " paragraph
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$p6.Style = "Normal"
$p6.Range.Text = "Normal text "

# 5c) blank paragraph
$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Style = "Normal"

# 5d) "Ordered List" paragraph, styled as List Bullet
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Item(8)
$p8.Range.Text = "Ordered List"
$p8.Style = "List Bullet"

# 5e) "Numbered List" paragraph, styled as List Number
$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Numbered List"
$p9.Style = "List Number"

# 6) Finally, remove the demo table itself.
$d.Tables.Item(1).Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
Write-Output ("Final content: [" + $d.Content.Text + "]")
